# Commit: "update ảnh cho 20 sản phẩm" (update images for 20 products)
#
# The author removed 13 stale product-image filenames
# (CuaHangTrangSuc\productsInfo\products2\*.png and \products3\*.png) that no
# longer correspond to an existing image on disk. In the workbook, the cells
# that used to reference those specific filenames now store the literal
# string "null" instead (mirroring the "null" placeholder already used
# elsewhere in columns pic2/pic3 for products with a missing picture).
#
# Every other shared-string index shuffles automatically once those 13
# strings drop out of the shared-strings table - that is handled by the
# engine when we simply (re)write the affected cells' values, so we only
# need to touch the cells whose *content* actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cellsToNull = @("F22", "F24", "E25", "F25", "E27", "F27", "E33", "F33", "F34", "E38", "F38", "E39", "F39")
foreach ($addr in $cellsToNull) {
    $ws.Range($addr).Value = "null"
}

# Reflect the author's final cursor position / selection in the sheet view.
$ws.Range("E42").Select()
